$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new range to be plain text so numeric-looking strings
# ("2", "75.67", "0.00", ...) are preserved verbatim instead of being
# coerced into numbers (which would also drop formatting like trailing
# zeros).
$rng = $ws.Range("A30:E39")
$rng.NumberFormat = "@"

# A lone "'" forces Excel to commit an *empty text* cell (rather than a
# truly blank/null one) for the rows whose SKU column is blank in the
# source data.
$emptyText = "'"

$data = @(
  @("DAWN",      "Jam - Raspberry Pure",           "2",  "75.67",  "151.34"),
  @("Lentz",     "Oats",                           "1",  "43.94",  "43.94"),
  @($emptyText,  "Mustard - Honey",                "2",  "0.00",   "0.00"),
  @("Palmer",    "Butter - Salted",                "1",  "0.00",   "0.00"),
  @($emptyText,  "Goat Cheese",                    "12", "147.04", "1764.48"),
  @("PERF",      "Vegan Egg",                      "2",  "99.59",  "199.18"),
  @("Casa",      "Nuts - Pine",                    "2",  "127.50", "255.00"),
  @($emptyText,  "Nuts - Walnut Halves & Pieces",  "1",  "3.08",   "3.08"),
  @($emptyText,  "Sugar - Light Brown",            "2",  "46.70",  "93.40"),
  @($emptyText,  "Flour - Millers Choice",         "2",  "0.00",   "0.00")
)

$row = 30
foreach ($line in $data) {
    $ws.Cells.Item($row, 1).Value = $line[0]
    $ws.Cells.Item($row, 2).Value = $line[1]
    $ws.Cells.Item($row, 3).Value = $line[2]
    $ws.Cells.Item($row, 4).Value = $line[3]
    $ws.Cells.Item($row, 5).Value = $line[4]
    $row++
}
